# Adding the changes we made on may 9th
#
# Prepends one newly-captured sensor reading at the top of the "struggle"
# data block (new row 2), which pushes the previously-recorded ax/ay/az/
# gx/gy/gz measurements down by one row, and appends nine more freshly
# captured rows at the bottom (timestamps 2000-2900). The timestamp/label
# columns keep running in sequence; only the six sensor-measurement columns
# (C:H) shift down for the pre-existing rows. Net effect: the sheet grows
# from 20 data rows (A1:H21) to 30 data rows (A1:H31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement block for A2:H31: timestamp, label, ax, ay, az, gx, gy, gz
$rows = @(
    @(0,"struggle",0.2804546356201172,0.4303635954856872,-0.691750168800354,0.1050096067542932,-1.756468223065746,0.4945203567645989),
    @(100,"struggle",0.1987819671630859,0.2879692316055298,-0.9282988905906676,-0.0286234012063665,-0.7998002785809195,0.0811297598541999),
    @(200,"struggle",0.2130537033081054,0.492556095123291,-0.9500083923339844,0.03028146071093407,0.6503314929349078,-0.2584614434412545),
    @(300,"struggle",5.370540142059326,-2.71519660949707,-10.77233219146728,0.6556547600395809,3.617322595752007,-1.325163067603595),
    @(400,"struggle",1.807073593139648,1.692888975143433,0.979395866394043,0.9173419645854408,4.483878631981052,-1.026156997194096),
    @(500,"struggle",-1.115023136138916,-0.1916569471359253,-1.226519584655761,0.6556111379545562,2.95567157317181,0.3850134197546516),
    @(600,"struggle",-0.6978027820587158,-0.2956833839416504,0.6510324478149414,0.01345462458474067,2.573593986277679,-0.1254206044333366),
    @(700,"struggle",-3.345348596572876,1.023304224014282,2.694005012512207,-0.0735905529284972,1.541839412736652,-0.5845232478209913),
    @(800,"struggle",-9.996002197265623,-0.0120857954025268,1.461801767349243,0.5347626531610679,-1.351723280640276,0.006906517306153076),
    @(900,"struggle",4.229874610900879,-2.616225481033325,-11.63714218139648,0.2076723618166787,-3.982405117579868,-0.1024508295314652),
    @(1000,"struggle",-5.585046291351318,-1.405134916305542,1.352597236633301,-0.6092913156869469,-4.47834036301593,0.3969252789203029),
    @(1100,"struggle",-3.109971046447754,-0.2804408073425293,-1.378412246704102,-0.4407143246154432,-3.636212723595754,0.6680403655889118),
    @(1200,"struggle",-6.216216564178467,2.080293655395508,-6.55400562286377,0.2482107579708095,-2.94579162889598,0.002957710806204883),
    @(1300,"struggle",-2.478143692016602,-0.5579037666320801,-6.963788509368896,0.1622283559064473,0.07419828006199625,-0.1648120685499541),
    @(1400,"struggle",1.977913856506348,1.696393609046936,4.733741760253906,0.2419649730531533,4.375038618944129,-0.3818312956362353),
    @(1500,"struggle",-1.786364078521728,-0.5252545475959778,-2.082462787628174,0.5893073264433398,6.328164937544838,0.07646724885823974),
    @(1600,"struggle",-3.566806316375732,0.4948284029960632,1.737120628356934,-0.175296502453943,2.29715876919883,1.448470601013731),
    @(1700,"struggle",-6.714832305908203,2.814615249633789,5.518161296844482,0.8341801327710269,0.753240480714915,-0.1038439748238527),
    @(1800,"struggle",-4.827847480773926,-7.360194683074951,7.510563373565674,-0.1397884144466732,0.2036113617371567,-0.8063265510967792),
    @(1900,"struggle",-1.930130958557129,4.015055179595947,-5.047464370727539,-0.7627712436476506,-1.678813515877241,-1.093177660387389),
    @(2000,"struggle",-4.561958312988281,0.1262733936309814,-1.993290901184082,-0.0004051567948546442,-2.697911312993704,-0.4782794093599184),
    @(2100,"struggle",-2.729169845581055,3.059413433074951,-4.533473014831543,-0.8633800915309378,-0.3154059344408438,0.4824308418497782),
    @(2200,"struggle",-0.4514303207397461,-0.07753515243530271,-1.056098580360413,-0.4081483519807154,-0.6726997543354425,-0.2190668820118418),
    @(2300,"struggle",1.037992477416992,-1.273390769958496,0.4362349510192871,0.2211332225373814,0.241335413285664,0.08368853798934378),
    @(2400,"struggle",0.0754270553588867,1.646718859672546,1.695090532302856,0.06768137718341787,0.3379019900244107,0.1505034766635118),
    @(2500,"struggle",-0.2560558319091797,0.3026316165924072,-0.4233262538909912,0.07254024853511698,0.5556785336562575,-0.05807583201296457),
    @(2600,"struggle",0.6335611343383789,0.8106564879417419,-1.443797469139099,0.1816357883567719,0.1322741392923868,-0.08515337003128903),
    @(2700,"struggle",0.09285736083984369,0.7357764840126038,-1.646607518196106,-0.02734556931013958,-0.1169588795425942,0.04497027853313797),
    @(2800,"struggle",0.0882749557495117,0.1726978719234466,-0.9354652166366576,-0.02540700723017953,-0.06986615411481072,-0.074921377335808),
    @(2900,"struggle",0.2656211853027344,0.4902379512786865,-0.8409426212310791,0.02237761537639455,-0.07008743807863513,-0.003453258577050004)
)

$nRows = $rows.Count
$nCols = 8
$arr = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $arr[$r,$c] = $rows[$r][$c]
    }
}

$ws.Range("A2:H31").Value = $arr

Write-Output "Wrote $nRows rows to A2:H31 (dimension now A1:H31)"
